$d = $word.ActiveDocument

# Locate the "Emotional_EMA" heading paragraph (by index, since Paragraph.Next
# does not yield a properly-anchored Range in this host) and then target the
# start of the paragraph right after it -- the old trailing empty paragraph --
# as the insertion point for the new OOXML content.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Emotional_EMA") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Emotional_EMA' heading paragraph"
}

$insertionPoint = $d.Paragraphs.Item($targetIndex + 1).Range
$insertionPoint.Collapse(1)

$newContentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:r>
        <w:t>From readme…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>‘’’</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>This Electromagnetic Articulography (EMA) database includes articulatory motions recorded by an EMA system.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Talkers produced simulated (acted) emotional speech.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>A set of 10 sentences was commonly used for speech recording of a male (AB) and two females (JN, LS), who are native speakers of American English.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>On top of the 10 sentences, there are 4 additional sentences used for recording by only AB.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Each sentence was produced five times for four different emotions, such as neutrality, anger, sadness and happiness.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">In </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>totol</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, AB produced 280 utterances (14 sentences x 5 repetitions x 4 emotions), and JR and JN produced 200 utterances (10 sentences x 5 repetitions x 4 emotions).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Each utterance was digitalized in 12-bit amplitude resolution with 16kHz sampling rate.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Speech was recorded simultaneously by the EMA system so that speech and corresponding articulatory movements are aligned in time.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>‘’’</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Created </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>valence_scores_per_sample</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> from DocumentationEma.txt</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>In DocumentationEma.txt, there were two filenames misspelled:</w:t>
      </w:r>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblW w:w="5185" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="1715"/>
        <w:gridCol w:w="960"/>
        <w:gridCol w:w="2831"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1580" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
              <w:left w:val="nil"/>
              <w:bottom w:val="nil"/>
              <w:right w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>4EMO_~43.WAV</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="960" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
              <w:left w:val="nil"/>
              <w:bottom w:val="nil"/>
              <w:right w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>=</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="2645" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
              <w:left w:val="nil"/>
              <w:bottom w:val="nil"/>
              <w:right w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>4emo_ls_angry_41_041.wav</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1580" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
              <w:left w:val="nil"/>
              <w:bottom w:val="nil"/>
              <w:right w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>4EMO_~86.WAV</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="960" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
              <w:left w:val="nil"/>
              <w:bottom w:val="nil"/>
              <w:right w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>=</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="2645" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
              <w:left w:val="nil"/>
              <w:bottom w:val="nil"/>
              <w:right w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>4emo_ls_happy_32_032.wav</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>For the emotion category votes, each sample obtained a majority vote for one category (3</w:t>
      </w:r>
      <w:r>
        <w:t>/4</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> or 4</w:t>
      </w:r>
      <w:r>
        <w:t>/4</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> votes)</w:t>
      </w:r>
      <w:r>
        <w:t>. All of these matched the intended emotion.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> These were from t</w:t>
      </w:r>
      <w:r>
        <w:t>he best_xxx_files.txt files</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Samples were also rated on valence separately (different evaluators).</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> These ratings didn’t always match </w:t>
      </w:r>
      <w:r>
        <w:t>the valence of the intended emotion.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> I kept the samples where either the majority valence vote (if present) or the average valence rating matched the valence of the intended emotion.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> I allowed both criteria to increase the number of samples retained.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> 32 samples were discarded because of perceived-intended mismatch.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The best_xxx_files.txt files do not contain all the files listed in </w:t>
      </w:r>
      <w:r>
        <w:t>DocumentationEma.txt</w:t>
      </w:r>
      <w:r>
        <w:t>!</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> The</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> leftovers</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> were assessed by valence only.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Discarded 58</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>EmoV-DB_sorted</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newContentXml) | Out-Null

Write-Output "Inserted new content after the 'Emotional_EMA' heading."
